# Rearranges the subject/period values in the "RESULT" timetable sheet.
# The grid layout (days / periods / classes) stays exactly the same; only
# the text values in certain cells are updated to their new subject.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Phys1-1"
$ws.Range("B2").Value = "Russian1-2"
$ws.Range("D2").Value = "English1-4"
$ws.Range("E2").Value = "Math1-5"

$ws.Range("A3").Value = "Russian1-1"
$ws.Range("B3").Value = "Math1-2"
$ws.Range("D3").Value = "Math1-4"
$ws.Range("E3").Value = "Litr1-5"

$ws.Range("A4").Value = "Math1-1"
$ws.Range("E4").Value = "Russian1-5"

$ws.Range("B7").Value = "Math2-2"
$ws.Range("C7").Value = "Math2-3"
$ws.Range("D7").Value = "Phys2-4"

$ws.Range("B8").Value = "Phys2-2"
$ws.Range("C8").Value = "Phys2-3"

$ws.Range("B11").Value = "Russian3-2"
$ws.Range("C11").Value = "Phys3-3"

$ws.Range("A12").Value = "English3-1"
$ws.Range("B12").Value = "Phys3-2"
$ws.Range("C12").Value = "Russian3-3"
$ws.Range("E12").Value = "English3-5"

$ws.Range("A13").Value = "Phys3-1"

$ws.Range("B16").Value = "English4-2"
$ws.Range("C16").Value = "Math4-3"
$ws.Range("D16").Value = "Phys4-4"

$ws.Range("C17").Value = "Phys4-3"
$ws.Range("D17").Value = "Math4-4"

$ws.Range("B20").Value = "Phys5-2"
$ws.Range("D20").Value = "Math5-4"
$ws.Range("E20").Value = "Phys5-5"
